# Applies the Week02 task #2 edit: reorders/updates the triples in the
# triple_sheet worksheet and appends 4 new rows (commander/officer triples).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired contents for A1:C14 (header row kept as-is).
$data = @(
    @("subject",        "predicate", "object"),
    @("peiper:Person",  "born_in",   "Wilmersdorf, Berlin, Germany"),
    @("himler:Person",  "know",      "peiper"),
    @("himler:Person",  "born",      "7 October 1900"),
    @("himler:Person",  "die",       "23 May 1945"),
    @("peiper:Person",  "is_a",      "soldier"),
    @("peiper:Person",  "is_a",      "nazi member"),
    @("peiper:Person",  "born",      "30 January 1915"),
    @("peiper:Person",  "is_a",      "German"),
    @("peiper:Person",  "die",       "14 July 1976"),
    @("himler:Person",  "is_a",      "officer"),
    @("himler:Person",  "is_a",      "commander"),
    @("peiper:Person",  "is_a",      "commander"),
    @("peiper:Person",  "is_a",      "officer")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
